$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "Sheet1" to "raw data"
$ws.Name = "raw data"

# Fix the header text in D1: "Answer_relevance" -> "Answer relevance"
$ws.Range("D1").Value = "Answer relevance"

# Move the active cell selection from D1 to D2
$ws.Range("D2").Select()
